{"js": "// Replace the outdated representative's name with the new one in the\n// body of the document (\"neste ato representada por ...\").\nconst body = context.document.body;\nconst searchResults = body.search(\"Jo\u00e3o Paulo Salazar Dias\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"Ana Sofia Rodrigues dos Reis Mota\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the outdated representative's name with the new one in the\n# body of the document (\"neste ato representada por ...\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# wdReplaceAll = 2, wdFindContinue = 1\n$find.Execute(\n    \"Jo\u00e3o Paulo Salazar Dias\",  # FindText\n    $true,                       # MatchCase\n    $false,                      # MatchWholeWord\n    $false,                      # MatchWildcards\n    $false,                      # MatchSoundsLike\n    $false,                      # MatchAllWordForms\n    $true,                       # Forward\n    1,                           # Wrap -> wdFindContinue\n    $false,                      # Format\n    \"Ana Sofia Rodrigues dos Reis Mota\",  # ReplaceWith\n    2                            # Replace -> wdReplaceAll\n) | Out-Null\n"}
